$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Il18"
$ws.Cells.Item(2, 3).Value = "Il18rap"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.855689333333333
$ws.Cells.Item(2, 8).Value = 11.567068
$ws.Cells.Item(2, 9).Value = 0.1340401150840085
$ws.Cells.Item(2, 10).Value = 0.1340401150840085
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1125586666666667
$ws.Cells.Item(2, 14).Value = 0.337676
$ws.Cells.Item(2, 15).Value = 0.08976149178944288
$ws.Cells.Item(2, 16).Value = 0.0897614917894429
$ws.Cells.Item(2, 17).Value = 0.4339912504408888
$ws.Cells.Item(2, 18).Value = 3.905921253967999
$ws.Cells.Item(2, 19).Value = 0.01203164068956921
$ws.Cells.Item(2, 20).Value = 0.01203164068956921

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Il18"
$ws.Cells.Item(3, 3).Value = "Il18rap"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.855689333333333
$ws.Cells.Item(3, 8).Value = 11.567068
$ws.Cells.Item(3, 9).Value = 0.1340401150840085
$ws.Cells.Item(3, 10).Value = 0.1340401150840085
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.3610756666666666
$ws.Cells.Item(3, 14).Value = 1.083227
$ws.Cells.Item(3, 15).Value = 0.287944868651023
$ws.Cells.Item(3, 16).Value = 0.287944868651023
$ws.Cells.Item(3, 17).Value = 1.392195596492888
$ws.Cells.Item(3, 18).Value = 12.529760368436
$ws.Cells.Item(3, 19).Value = 0.03859616333183284
$ws.Cells.Item(3, 20).Value = 0.03859616333183285

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Il18"
$ws.Cells.Item(4, 3).Value = "Il18rap"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.855689333333333
$ws.Cells.Item(4, 8).Value = 11.567068
$ws.Cells.Item(4, 9).Value = 0.1340401150840085
$ws.Cells.Item(4, 10).Value = 0.1340401150840085
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.7803406666666667
$ws.Cells.Item(4, 14).Value = 2.341022
$ws.Cells.Item(4, 15).Value = 0.622293639559534
$ws.Cells.Item(4, 16).Value = 0.622293639559534
$ws.Cells.Item(4, 17).Value = 3.008751184832889
$ws.Cells.Item(4, 18).Value = 27.078760663496
$ws.Cells.Item(4, 19).Value = 0.08341231106260645
$ws.Cells.Item(4, 20).Value = 0.08341231106260646

# Row 5
$ws.Cells.Item(5, 1).Value = "M2"
$ws.Cells.Item(5, 2).Value = "Il18"
$ws.Cells.Item(5, 3).Value = "Il18rap"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 22.00607133333333
$ws.Cells.Item(5, 8).Value = 66.018214
$ws.Cells.Item(5, 9).Value = 0.7650243780187601
$ws.Cells.Item(5, 10).Value = 0.7650243780187601
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.1125586666666667
$ws.Cells.Item(5, 14).Value = 0.337676
$ws.Cells.Item(5, 15).Value = 0.08976149178944288
$ws.Cells.Item(5, 16).Value = 0.0897614917894429
$ws.Cells.Item(5, 17).Value = 2.476974047851555
$ws.Cells.Item(5, 18).Value = 22.292766430664
$ws.Cells.Item(5, 19).Value = 0.06866972942625459
$ws.Cells.Item(5, 20).Value = 0.0686697294262546

# Row 6
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Il18"
$ws.Cells.Item(6, 3).Value = "Il18rap"
$ws.Cells.Item(6, 4).Value = "M2"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 22.00607133333333
$ws.Cells.Item(6, 8).Value = 66.018214
$ws.Cells.Item(6, 9).Value = 0.7650243780187601
$ws.Cells.Item(6, 10).Value = 0.7650243780187601
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.3610756666666666
$ws.Cells.Item(6, 14).Value = 1.083227
$ws.Cells.Item(6, 15).Value = 0.287944868651023
$ws.Cells.Item(6, 16).Value = 0.287944868651023
$ws.Cells.Item(6, 17).Value = 7.945856877397556
$ws.Cells.Item(6, 18).Value = 71.512711896578
$ws.Cells.Item(6, 19).Value = 0.2202848440434425
$ws.Cells.Item(6, 20).Value = 0.2202848440434425

# Row 7
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Il18"
$ws.Cells.Item(7, 3).Value = "Il18rap"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 22.00607133333333
$ws.Cells.Item(7, 8).Value = 66.018214
$ws.Cells.Item(7, 9).Value = 0.7650243780187601
$ws.Cells.Item(7, 10).Value = 0.7650243780187601
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.7803406666666667
$ws.Cells.Item(7, 14).Value = 2.341022
$ws.Cells.Item(7, 15).Value = 0.622293639559534
$ws.Cells.Item(7, 16).Value = 0.622293639559534
$ws.Cells.Item(7, 17).Value = 17.17223237496756
$ws.Cells.Item(7, 18).Value = 154.550091374708
$ws.Cells.Item(7, 19).Value = 0.476069804549063
$ws.Cells.Item(7, 20).Value = 0.476069804549063

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Il18"
$ws.Cells.Item(8, 3).Value = "Il18rap"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.903429
$ws.Cells.Item(8, 8).Value = 8.710287000000001
$ws.Cells.Item(8, 9).Value = 0.1009355068972313
$ws.Cells.Item(8, 10).Value = 0.1009355068972313
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.1125586666666667
$ws.Cells.Item(8, 14).Value = 0.337676
$ws.Cells.Item(8, 15).Value = 0.08976149178944288
$ws.Cells.Item(8, 16).Value = 0.0897614917894429
$ws.Cells.Item(8, 17).Value = 0.3268060970013333
$ws.Cells.Item(8, 18).Value = 2.941254873012
$ws.Cells.Item(8, 19).Value = 0.009060121673619083
$ws.Cells.Item(8, 20).Value = 0.009060121673619085

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Il18"
$ws.Cells.Item(9, 3).Value = "Il18rap"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.903429
$ws.Cells.Item(9, 8).Value = 8.710287000000001
$ws.Cells.Item(9, 9).Value = 0.1009355068972313
$ws.Cells.Item(9, 10).Value = 0.1009355068972313
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.3610756666666666
$ws.Cells.Item(9, 14).Value = 1.083227
$ws.Cells.Item(9, 15).Value = 0.287944868651023
$ws.Cells.Item(9, 16).Value = 0.287944868651023
$ws.Cells.Item(9, 17).Value = 1.048357561794333
$ws.Cells.Item(9, 18).Value = 9.435218056149
$ws.Cells.Item(9, 19).Value = 0.0290638612757477
$ws.Cells.Item(9, 20).Value = 0.0290638612757477

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Il18"
$ws.Cells.Item(10, 3).Value = "Il18rap"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.903429
$ws.Cells.Item(10, 8).Value = 8.710287000000001
$ws.Cells.Item(10, 9).Value = 0.1009355068972313
$ws.Cells.Item(10, 10).Value = 0.1009355068972313
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.7803406666666667
$ws.Cells.Item(10, 14).Value = 2.341022
$ws.Cells.Item(10, 15).Value = 0.622293639559534
$ws.Cells.Item(10, 16).Value = 0.622293639559534
$ws.Cells.Item(10, 17).Value = 2.265663721479334
$ws.Cells.Item(10, 18).Value = 20.390973493314
$ws.Cells.Item(10, 19).Value = 0.06281152394786452
$ws.Cells.Item(10, 20).Value = 0.06281152394786452
